$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Update the EmailMultiLocation "ReSchedule CampName" value (J2)
#    MultiLocation Test06/24/2020 11:06:56 AM
#      -> MultiLocation Test06/25/2020 7:05:30 AM
# ------------------------------------------------------------------
$eml = $wb.Worksheets.Item("EmailMultiLocation")
$eml.Range("J2").Value = "MultiLocation Test06/25/2020 7:05:30 AM"

# ------------------------------------------------------------------
# 2. Insert a new "GRLMultiLocation" sheet right after GRLBrandCampaign
# ------------------------------------------------------------------
$afterGrl = $wb.Worksheets.Item("GRLBrandCampaign")
$grlMulti = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterGrl)
$grlMulti.Name = "GRLMultiLocation"

$grlMulti.Range("A1").Value = "CamType"
$grlMulti.Range("B1").Value = "CamOption"
$grlMulti.Range("C1").Value = "CamLang"
$grlMulti.Range("D1").Value = "CamName"
$grlMulti.Range("E1").Value = "CamDes"
$grlMulti.Range("F1").Value = "Locations"
$grlMulti.Range("G1").Value = "1Star Messaging"
$grlMulti.Range("H1").Value = "3Star Messaging"
$grlMulti.Range("I1").Value = "4Star Messaging"
$grlMulti.Range("J1").Value = "ProcessedCampaign Name"
$grlMulti.Range("A1:J1").Font.Bold = $true

$grlMulti.Range("A2").Value = "General Review Link"
$grlMulti.Range("B2").Value = "Multi-location"
$grlMulti.Range("C2").Value = "English"
$grlMulti.Range("D2").Value = "GRL MultiLocation Test"
$grlMulti.Range("E2").Value = "Test"
$grlMulti.Range("F2").Value = "9000334772,9000334773,9000334774"
$grlMulti.Range("G2").Value = "Sorry for your experience and we'll make sure you will have better experience next time"
$grlMulti.Range("H2").Value = "Thanks you for feedback"
$grlMulti.Range("I2").Value = "Thanks for your feedback! Have a greatday"
$grlMulti.Range("J2").Value = "GRL MultiLocation Test06/25/2020 8:05:25 AM"

for ($i = 1; $i -le 10; $i++) {
    $grlMulti.Columns.Item($i).AutoFit()
}

# ------------------------------------------------------------------
# 3. Append a new "Filters" sheet at the end of the workbook
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$filters = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$filters.Name = "Filters"

$filters.Range("A1").Value = "FilterType"
$filters.Range("B1").Value = "FilterSetUp"
$filters.Range("A1:B1").Font.Bold = $true

$filters.Range("A2").Value = "Email"
$filters.Range("B2").Value = "null"
$filters.Range("A3").Value = "Email"
$filters.Range("B3").Value = "Multi-location"
$filters.Range("A4").Value = "General Review Link"
$filters.Range("B4").Value = "null"
$filters.Range("A5").Value = "General Review Link"
$filters.Range("B5").Value = "Multi-location"

$filters.Columns.Item(1).AutoFit()
$filters.Columns.Item(2).AutoFit()

$filters.Range("A6:XFD14").Select()

# ------------------------------------------------------------------
# 4. Workbook view: active sheet = Filters
# ------------------------------------------------------------------
$filters.Select()
